$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append S30 sprint task rows 353-371 ---
# Groups: G08 'observability & performance (sec 12)', G09 'restart & operational
# resilience (sec 13)', G10 'pause & resume semantics (sec 14)'.

$rowCount = 19
$startRow = 353

# Columns A-E (sprint#, group#, group task description, task#, task description)
$mainData = New-Object 'object[,]' 19,5

# Column G (status)
$statusData = New-Object 'object[,]' 19,1

$mainData[0,0] = 'S30'
$mainData[0,1] = 'G08'
$mainData[0,2] = 'Strategy deployment v3 — observability & performance (sec 12)'
$mainData[0,3] = 'S30_G08_TB001'
$mainData[0,4] = 'Backend: Add deployment heartbeat fields in DB (last_eval_at, last_eval_bar_end_ts, runtime_state enum FLAT/IN_POSITION/WARMING_UP/PAUSED_*/ERROR, last_decision enum, last_decision_reason, next_eval_at) + expose in API schemas; index last_eval_at/next_eval_at for table views.'
$statusData[0,0] = 'pending'

$mainData[1,0] = 'S30'
$mainData[1,1] = 'G08'
$mainData[1,2] = 'Strategy deployment v3 — observability & performance (sec 12)'
$mainData[1,3] = 'S30_G08_TB002'
$mainData[1,4] = 'Backend: Add per-deployment append-only event journal table (deployment_event_log) + writer helpers; emit events for BAR_CLOSED received, eval start/finish, entry/exit signal booleans, order intent + submit/fill states, risk exits, reconciliation actions.'
$statusData[1,0] = 'pending'

$mainData[2,0] = 'S30'
$mainData[2,1] = 'G08'
$mainData[2,2] = 'Strategy deployment v3 — observability & performance (sec 12)'
$mainData[2,3] = 'S30_G08_TB003'
$mainData[2,4] = 'Backend: Ensure order intent → order → fills/trades are linked by deployment_id + intent_id + dedupe_key (paper and live); add missing columns/migrations and propagate through execution paths.'
$statusData[2,0] = 'pending'

$mainData[3,0] = 'S30'
$mainData[3,1] = 'G08'
$mainData[3,2] = 'Strategy deployment v3 — observability & performance (sec 12)'
$mainData[3,3] = 'S30_G08_TB004'
$mainData[3,4] = 'Backend: Implement live per-deployment performance summary service + API (realized/unrealized P&L, current position, trade count, last trade time; optional drawdown) computed from fills/positions and cached safely.'
$statusData[3,0] = 'pending'

$mainData[4,0] = 'S30'
$mainData[4,1] = 'G08'
$mainData[4,2] = 'Strategy deployment v3 — observability & performance (sec 12)'
$mainData[4,3] = 'S30_G08_TB005'
$mainData[4,4] = 'Backend: Add deployment equity curve points table + API; write points on each exit and on EOD/proxy close (mark-to-market) for RUNNING deployments.'
$statusData[4,0] = 'pending'

$mainData[5,0] = 'S30'
$mainData[5,1] = 'G08'
$mainData[5,2] = 'Strategy deployment v3 — observability & performance (sec 12)'
$mainData[5,3] = 'S30_G08_TF001'
$mainData[5,4] = 'Frontend: Deployments table shows State, Last Eval, Last Decision (chips) + tooltips/help icons explaining state/decision; refresh cadence and “evidence of life” UX.'
$statusData[5,0] = 'pending'

$mainData[6,0] = 'S30'
$mainData[6,1] = 'G08'
$mainData[6,2] = 'Strategy deployment v3 — observability & performance (sec 12)'
$mainData[6,3] = 'S30_G08_TF002'
$mainData[6,4] = 'Frontend: Add deployment details right-side drawer (recommended) with tabs Summary/Equity/Journal/Orders/Trades/Diagnostics; ensure selection from list opens drawer and deep-linking still works.'
$statusData[6,0] = 'pending'

$mainData[7,0] = 'S30'
$mainData[7,1] = 'G08'
$mainData[7,2] = 'Strategy deployment v3 — observability & performance (sec 12)'
$mainData[7,3] = 'S30_G08_TD001'
$mainData[7,4] = 'Docs: Update strategy deployment help to explain heartbeat fields, event journal, P&L attribution, and equity curve semantics (paper vs live caveats).'
$statusData[7,0] = 'pending'

$mainData[8,0] = 'S30'
$mainData[8,1] = 'G08'
$mainData[8,2] = 'Strategy deployment v3 — observability & performance (sec 12)'
$mainData[8,3] = 'S30_G08_TT001'
$mainData[8,4] = 'Tests: Heartbeat updates per eval (including NO_BAR/MARKET_CLOSED decisions), event journal insertion, order/fill linkage, and performance/equity API correctness (paper path).'
$statusData[8,0] = 'pending'

$mainData[9,0] = 'S30'
$mainData[9,1] = 'G09'
$mainData[9,2] = 'Strategy deployment v3 — restart & operational resilience (sec 13)'
$mainData[9,3] = 'S30_G09_TB001'
$mainData[9,4] = 'Backend: Implement reconciliation-on-start routine: on BE startup load RUNNING deployments, fetch broker/paper positions+open orders, compare expected vs actual, repair safe mismatches and PAUSE deployment when unsafe; record journal events.'
$statusData[9,0] = 'pending'

$mainData[10,0] = 'S30'
$mainData[10,1] = 'G09'
$mainData[10,2] = 'Strategy deployment v3 — restart & operational resilience (sec 13)'
$mainData[10,3] = 'S30_G09_TB002'
$mainData[10,4] = 'Backend: Audit and harden idempotency/dedupe across all trading actions (evaluation, order intent create, submission, trailing updates, MIS auto-flatten); unify deterministic dedupe keys and add guardrails to prevent duplicates on restart.'
$statusData[10,0] = 'pending'

$mainData[11,0] = 'S30'
$mainData[11,1] = 'G09'
$mainData[11,2] = 'Strategy deployment v3 — restart & operational resilience (sec 13)'
$mainData[11,3] = 'S30_G09_TB003'
$mainData[11,4] = 'Backend: Add operational controls/telemetry for restarts (e.g., reconcile-only startup mode flag, last_reconcile_at/result on health endpoint, and rate-limited reconciliation sweeps outside market hours).'
$statusData[11,0] = 'pending'

$mainData[12,0] = 'S30'
$mainData[12,1] = 'G09'
$mainData[12,2] = 'Strategy deployment v3 — restart & operational resilience (sec 13)'
$mainData[12,3] = 'S30_G09_TD001'
$mainData[12,4] = 'Docs: Operator runbook for restart safety (off-hours guarantees, recommended pause→restart→resume workflow, and how reconciliation resolves/pauses deployments).'
$statusData[12,0] = 'pending'

$mainData[13,0] = 'S30'
$mainData[13,1] = 'G09'
$mainData[13,2] = 'Strategy deployment v3 — restart & operational resilience (sec 13)'
$mainData[13,3] = 'S30_G09_TT001'
$mainData[13,4] = 'Tests: Restart simulation (persist jobs/actions/orders then re-init workers) verifies no duplicate orders, reconciliation runs for RUNNING deployments, and unsafe mismatches lead to PAUSED + journal entry.'
$statusData[13,0] = 'pending'

$mainData[14,0] = 'S30'
$mainData[14,1] = 'G10'
$mainData[14,2] = 'Strategy deployment v3 — pause & resume semantics (sec 14)'
$mainData[14,3] = 'S30_G10_TB001'
$mainData[14,4] = 'Backend: Add explicit Pause/Resume endpoints and persist paused_at/resumed_at + reason; ensure status transitions are audited and reflected in heartbeat fields.'
$statusData[14,0] = 'pending'

$mainData[15,0] = 'S30'
$mainData[15,1] = 'G10'
$mainData[15,2] = 'Strategy deployment v3 — pause & resume semantics (sec 14)'
$mainData[15,3] = 'S30_G10_TB002'
$mainData[15,4] = 'Backend: Enforce scheduler/worker double-gate: scheduler must not enqueue BAR_CLOSED/DAILY_PROXY jobs for PAUSED deployments; worker must re-check status before execute; run-now is blocked (or returns “paused”).'
$statusData[15,0] = 'pending'

$mainData[16,0] = 'S30'
$mainData[16,1] = 'G10'
$mainData[16,2] = 'Strategy deployment v3 — pause & resume semantics (sec 14)'
$mainData[16,3] = 'S30_G10_TB003'
$mainData[16,4] = 'Backend: Ensure PAUSE does not cancel broker-side protections (disaster SL/GTT) and does not disable MIS auto-flatten schedule; document and test invariants.'
$statusData[16,0] = 'pending'

$mainData[17,0] = 'S30'
$mainData[17,1] = 'G10'
$mainData[17,2] = 'Strategy deployment v3 — pause & resume semantics (sec 14)'
$mainData[17,3] = 'S30_G10_TF001'
$mainData[17,4] = 'Frontend: Add Pause/Resume/Stop actions in deployments list and details; show paused-at timestamp and “protections remain active” notice; add help icon clarifying pause semantics.'
$statusData[17,0] = 'pending'

$mainData[18,0] = 'S30'
$mainData[18,1] = 'G10'
$mainData[18,2] = 'Strategy deployment v3 — pause & resume semantics (sec 14)'
$mainData[18,3] = 'S30_G10_TT001'
$mainData[18,4] = 'Tests: Paused deployments do not generate/execute BAR_CLOSED; run-now respects pause; protections are not removed; UI/API permission checks.'
$statusData[18,0] = 'pending'

$endRow = $startRow + $rowCount - 1
$mainRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 5))
$mainRange.Value2 = $mainData

$statusRange = $ws.Range($ws.Cells.Item($startRow, 7), $ws.Cells.Item($endRow, 7))
$statusRange.Value2 = $statusData

# Match the plain/default formatting used by the surrounding rows (no wrap style).
# Only touch columns A-E and G so we don't create a stray empty styled cell in F.
$mainRange.Style = "Normal"
$statusRange.Style = "Normal"
